# Fixes species-name labels (stripping stray trailing underscores/dashes
# left over from OCR'd source text), fills in the three previously-blank
# species rows (Rockfish, Salmon, Tuna bluefin), and updates the sheet's
# view/selection state, per the "party vessels catch script and minor
# changes to some xlsx" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clean up trailing punctuation/underscores in existing species labels ---
$ws.Range("A4").Value  = "Bofuto, Pacific"
$ws.Range("A5").Value  = "Cabezon"
$ws.Range("A6").Value  = "Halibut California"
$ws.Range("A7").Value  = "Lingcod"
$ws.Range("A13").Value = "Seabass white"
$ws.Range("A18").Value = "Yefiowtail"

# --- Fill in the previously blank species rows ---
$ws.Range("A10").Value = "Rockfish"
$ws.Range("A11").Value = "Salmon"
$ws.Range("A16").Value = "Tuna, bluefin"

# --- Remaining label clean ups ---
$ws.Range("A19").Value = "All others"
$ws.Range("A23").Value = "Total number of angler hours"

# --- Update the sheet view / current selection ---
$ws.Range("A24").Select()
